# "Minor touch to image" - re-case the two "script" substrings inside the
# <script>alert("XSS")</script> sample text to "ScRiPt" on slide 2, splitting
# the single run into five runs (this mirrors what PowerPoint does when a
# user selectively retypes part of a run: the run is split around the
# edited span while the untouched spans keep their original run
# properties).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# The textbox lives 3 groups deep: Group 201 -> (flattened) GroupItems.
# GroupItems flattens every descendant shape of the group, so we can reach
# "TextBox 183" directly from the top-level group shape.
$grp = $s.Shapes.Item(6)
$box = $grp.GroupItems.Item(12)

$tr = $box.TextFrame.TextRange

# The box uses spAutoFit; remember its current (correct) height so we can
# restore it after splitting the text into runs below, since re-running
# layout while slicing the range nudges the autofit height.
$origHeight = $box.Height

# Sanity check / reset to the known starting text (curly quotes, as in the
# original run) before slicing it up.
$tr.Text = "<script>alert(" + [char]8220 + "XSS" + [char]8221 + ")</script>"

# Re-case the opening "script" (chars 2-7) -> "ScRiPt".
$tr.Characters(2, 6).Text = "ScRiPt"

# Re-case the closing "script" (chars 23-28, unchanged since the first
# edit didn't change the overall length) -> "ScRiPt".
$tr.Characters(23, 6).Text = "ScRiPt"

# Restore the shape's autofit height to what it was before the run split.
$box.Height = $origHeight
